$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")
$ws.Range("B124").Value = "test"
